# Update calendario de examenes sheet:
#  - Fill in previously "Pendiente" (pending) dates for row 8 (entrega de calificaciones)
#  - Fill in exam schedule details for several subjects (rows 9, 12-16)
#  - Swap the subject names/content of rows 10 and 11 (the "METODOLOGIA DE LA
#    PROGRAMACION" class replaces the old "DESARROLLO DE HABILIDADES DEL
#    PENSAMIENTO LOGICO" class in row 10, which now moves to row 11), per the
#    commit message "elimine una clase que no servia" (removed a class that
#    wasn't useful, effectively swapping which subject occupies which row and
#    updating schedules).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - Entrega de Calificaciones dates
$ws.Range("C8").Value = "31-Oct"
$ws.Range("D8").Value = "29-Nov"
$ws.Range("E8").Value = "4-Dic"
$ws.Range("F8").Value = "12-Dic"
$ws.Range("G8").Value = "17-Dic"

# Row 9 - ALGEBRA LINEAL
$ws.Range("C9").Value = "Miercoles 30 de Octubre `n Aula C-101 - 7:00"
$ws.Range("D9").Value = "Miercoles 27 de Noviembre `n Aula C-101 - 7:00"
$ws.Range("E9").Value = "Viernes 29 de Noviembre `n Aula C-101 - 7:00"
$ws.Range("F9").Value = "Martes 10 de Diciembre `n Aula C-101 - 7:00"
$ws.Range("G9").Value = "Viernes 13 de Diciembre `n Lab 5 - 12:40"

# Row 10 - now METODOLOGIA DE LA PROGRAMACION (was DESARROLLO DE HABILIDADES DEL PENSAMIENTO LOGICO)
$ws.Range("A10").Value = "METODOLOGÍA DE LA PROGRAMACIÓN"
$ws.Range("B10").Value = "Jueves 26 de Septiembre `n Lab 2 - 12:40"
$ws.Range("C10").Value = "Viernes 25 de Octubre `n Lab 6 - 9:00"
$ws.Range("D10").Value = "Jueves 28 de Noviembre `n Lab 2 - 12:40"
$ws.Range("E10").Value = "Jueves 28 de Noviembre `n Lab 2 - 12:40"
$ws.Range("F10").Value = "Miercoles 11 de Diciembre `n Lab 6 - 9:50"
$ws.Range("G10").Value = "Viernes 13 de Diciembre `n Lab 5 - 11:00"

# Row 11 - now DESARROLLO DE HABILIDADES DEL PENSAMIENTO LOGICO (was METODOLOGIA DE LA PROGRAMACION)
$ws.Range("A11").Value = "DESARROLLO DE HABILIDADES DEL PENSAMIENTO LÓGICO"
$ws.Range("B11").Value = "Viernes 27 de Septiembre `n Aula C-103 - 13:30"
$ws.Range("C11").Value = "Martes 29 de Octubre `n Aula C-101 - 12:40"
$ws.Range("D11").Value = "Martes 26 de Noviembre `n Aula C-101 - 12:40"
$ws.Range("E11").Value = "Martes 3 de Diciembre `n Aula C-101 - 12:40"
$ws.Range("F11").Value = "Miercoles 11 de Diciembre `n Aula C-101 - 7:50"
$ws.Range("G11").Value = "Viernes 13 de Diciembre `n Aula C-101 - 7:50"

# Row 12 - EXPRESION ORAL Y ESCRITA I
$ws.Range("C12").Value = "Lunes 28 de Octubre `n Aula C-101 - 7:00"
$ws.Range("D12").Value = "Lunes 25 de Noviembre `n Aula C-101 - 7:00"
$ws.Range("E12").Value = "Lunes 2 de Diciembre `n Aula C-101 - 7:00"
$ws.Range("F12").Value = "Miercoles 11 de Diciembre `n Lab 6 - 9:00"
$ws.Range("G12").Value = "Viernes 13 de Diciembre `n Lab 6 - 9:50"

# Row 13 - FUNDAMENTOS DE TI
$ws.Range("C13").Value = "Jueves 24 de Octubre `n Lab 2 - 9:00"
$ws.Range("D13").Value = "Miercoles 27 de Noviembre `n Lab 6 - 9:00"
$ws.Range("E13").Value = "Miercoles 27 de Noviembre `n Lab 6 - 9:00"
$ws.Range("F13").Value = "Miercoles 11 de Diciembre `n Lab1 - 11:00"
$ws.Range("G13").Value = "Viernes 13 de Diciembre `n Lab 6 - 9:00"

# Row 14 - FORMACION SOCIOCULTURAL I
$ws.Range("C14").Value = "Lunes 28 de Octubre `n Aula C-101 - 11:00"
$ws.Range("D14").Value = "Lunes 25 de Noviembre `n Aula C-101 - 11:00"
$ws.Range("E14").Value = "Lunes 2 de Diciembre `n Aula C-101 - 11:00"
$ws.Range("F14").Value = "Miercoles 11 de Diciembre `n Aula C-101 - 7:00"
$ws.Range("G14").Value = "Viernes 13 de Diciembre `n Aula C-101 - 7:00"

# Row 15 - FUNDAMENTOS DE REDES
$ws.Range("C15").Value = "Miercoles 30 de Octubre `n Lab1 - 11:00"
$ws.Range("D15").Value = "Viernes 22 de Noviembre `n Lab 5 - 11:00"
$ws.Range("E15").Value = "Viernes 29 de Noviembre `n Lab 5 - 11:00"
$ws.Range("F15").Value = "Miercoles 11 de Diciembre `n Lab1 - 11:50"
$ws.Range("G15").Value = "Viernes 13 de Diciembre `n Lab 5 - 11:50"

# Row 16 - INGLES I
$ws.Range("C16").Value = "De acuerdo a la coordinación de inglés"
$ws.Range("D16").Value = "De acuerdo a la coordinación de inglés"
$ws.Range("E16").Value = "De acuerdo a la coordinación de inglés"
$ws.Range("F16").Value = "De acuerdo a la coordinación de inglés"
$ws.Range("G16").Value = "De acuerdo a la coordinación de inglés"
